# "Generate Report for Handoff" - refresh the localization-status report:
#   - Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     642b4f16-960d-4aa9-ac0a-7fa31a23bbb2.md row (rows 7,8,10,11,12,14
#     share that same generated timestamp).
#   - zh-cn / de-de sheets: bump the matching "Latest Handoff Datetime"
#     for those same rows, and set the "Priority" for those rows to "ht"
#     (handoff type), now that they've gone through handoff generation.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 12, 14)

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-15 14:21:38"
}

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-15 14:21:30"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-15 14:21:38"
    $wsDeDe.Range("E$r").Value = "ht"
}
